$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 1.555537518143943
$ws.Range("E2").Value = 2.195459023351666
$ws.Range("F2").Value = 2.830872940803188
$ws.Range("G2").Value = 3.42534052304346
$ws.Range("H2").Value = 3.960342010682696
$ws.Range("I2").Value = 4.428373207691445
$ws.Range("J2").Value = 4.828127030019634
$ws.Range("K2").Value = 5.161176615309575
$ws.Range("L2").Value = 5.430169065316736
$ws.Range("M2").Value = 5.630231021774317
$ws.Range("N2").Value = 5.765338505637525
$ws.Range("O2").Value = 5.837136374553366
$ws.Range("P2").Value = 5.844757884475448
$ws.Range("Q2").Value = 5.800414334419341
$ws.Range("R2").Value = 5.727062937860904
$ws.Range("S2").Value = 5.639743751716479
$ws.Range("T2").Value = 5.548107449908807
$ws.Range("U2").Value = 5.458148760363308
$ws.Range("V2").Value = 5.373398162770581
$ws.Range("W2").Value = 5.295742448566643
$ws.Range("X2").Value = 5.225989775511186
$ws.Range("Y2").Value = 5.164257987690945
$ws.Range("Z2").Value = 5.110240195786897
$ws.Range("AA2").Value = 5.063384884311838
$ws.Range("AB2").Value = 5.023016443674727
$ws.Range("AC2").Value = 4.988414235200187
$ws.Range("AD2").Value = 4.958862911715991
$ws.Range("AE2").Value = 4.933682959501493
$ws.Range("AF2").Value = 4.915991456401172
